$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text formatting (avoid Excel
# auto-converting numeric-looking strings like "1.00" or "8.49" into numbers).
$cells = @('D2', 'E2', 'D3', 'E3', 'E4', 'D5', 'E5', 'D6', 'E6', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'E11', 'E12', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'E17', 'D18', 'E18', 'D19', 'E19', 'D20', 'E20', 'D21', 'E21', 'D22', 'E22', 'E23', 'D24', 'E24', 'D25', 'D26', 'E26', 'B27', 'C27', 'D27', 'E27', 'B28', 'C28', 'D28', 'E28', 'E29', 'E30', 'D31', 'E31', 'E32', 'D33', 'E33', 'E34', 'D35', 'E35', 'D36', 'E36', 'D37', 'E37', 'D38', 'E38', 'E39', 'D40', 'E40', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47', 'D48', 'E48', 'D49', 'E49', 'D50', 'E50', 'E51')
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.549.63'
$ws.Range('E2').Value = '  +2.31%  '
$ws.Range('D3').Value = '2.557.32'
$ws.Range('E3').Value = '  +5.07%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '572.42'
$ws.Range('E5').Value = '  +2.71%  '
$ws.Range('D6').Value = '150.29'
$ws.Range('E6').Value = '  +8.40%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  +0.72%  '
$ws.Range('D9').Value = '2.554.45'
$ws.Range('E9').Value = '  +4.95%  '
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('E13').Value = '  +3.77%  '
$ws.Range('D14').Value = '28.19'
$ws.Range('E14').Value = '  +9.39%  '
$ws.Range('D15').Value = '3.010.59'
$ws.Range('E15').Value = '  +5.18%  '
$ws.Range('D16').Value = '63.475.25'
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('E17').Value = '  +3.07%  '
$ws.Range('D18').Value = '2.551.82'
$ws.Range('E18').Value = '  +5.13%  '
$ws.Range('D19').Value = '11.69'
$ws.Range('E19').Value = '  +5.04%  '
$ws.Range('D20').Value = '342.41'
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').Value = '4.38'
$ws.Range('E21').Value = '  +3.45%  '
$ws.Range('D22').Value = '6.92'
$ws.Range('E22').Value = '  +1.22%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '66.19'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').Value = '0.171'
$ws.Range('D26').Value = '1.59'
$ws.Range('E26').Value = '  +4.26%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '8.49'
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  +8.06%  '
$ws.Range('E30').Value = '  +14.25%  '
$ws.Range('D31').Value = '0.0₃0840'
$ws.Range('E31').Value = '  +6.55%  '
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('D33').Value = '177.77'
$ws.Range('E33').Value = '  +3.86%  '
$ws.Range('E34').Value = '  +8.65%  '
$ws.Range('D35').Value = '415.54'
$ws.Range('E35').Value = '  +11.07%  '
$ws.Range('D36').Value = '0.408'
$ws.Range('E36').Value = '  +3.07%  '
$ws.Range('D37').Value = '19.14'
$ws.Range('E37').Value = '  +3.08%  '
$ws.Range('D38').Value = '4.46'
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '1.76'
$ws.Range('E40').Value = '  +4.28%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '40.02'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('D43').Value = '155.48'
$ws.Range('E43').Value = '  +6.45%  '
$ws.Range('D44').Value = '3.81'
$ws.Range('E44').Value = '  +3.77%  '
$ws.Range('D45').Value = '21.20'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('E46').Value = '  +3.79%  '
$ws.Range('D47').Value = '0.0534'
$ws.Range('E47').Value = '  +2.57%  '
$ws.Range('D48').Value = '0.0968'
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').Value = '0.0234'
$ws.Range('E49').Value = '  +5.53%  '
$ws.Range('D50').Value = '18.81'
$ws.Range('E50').Value = '  +4.60%  '
$ws.Range('E51').Value = '  +8.34%  '
